# Updates the cryptos price/volume table (GitHub Actions scheduled refresh).
# D (Price) and E (Volume 1h) columns hold values stored as text in the
# workbook, even when the text looks numeric (e.g. "243.92"). Plain
# `.Value = "243.92"` would let Excel auto-convert that into a real number
# cell, which would not match the source data (still plain text). For any
# new value that parses as a number we force the cell to remain text by
# flipping its number format to "@" before the assignment, then restore the
# original ("Normal") style afterwards so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "29.506.58"
$ws.Cells.Item(2, 5).Value = "  +1.76%  "
$ws.Cells.Item(3, 4).Value = "1.840.62"
$ws.Cells.Item(3, 5).Value = "  +0.61%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9977"
$ws.Cells.Item(4, 5).Value = "  -0.17%  "
Set-TextValue $ws.Cells.Item(5, 4) "243.92"
$ws.Cells.Item(5, 5).Value = "  +0.97%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.6275"
$ws.Cells.Item(6, 5).Value = "  +1.50%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.9987"
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.07423"
$ws.Cells.Item(8, 5).Value = "  -0.03%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.2948"
$ws.Cells.Item(9, 5).Value = "  +1.05%  "
Set-TextValue $ws.Cells.Item(10, 4) "23.48"
$ws.Cells.Item(10, 5).Value = "  +1.90%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07669"
$ws.Cells.Item(11, 5).Value = "  +0.21%  "
$ws.Cells.Item(12, 4).Value = "1.844.88"
$ws.Cells.Item(12, 5).Value = "  +1.08%  "
Set-TextValue $ws.Cells.Item(13, 4) "5.026"
$ws.Cells.Item(13, 5).Value = "  +0.70%  "
Set-TextValue $ws.Cells.Item(14, 4) "0.6822"
$ws.Cells.Item(14, 5).Value = "  +1.46%  "
Set-TextValue $ws.Cells.Item(15, 4) "83.48"
$ws.Cells.Item(15, 5).Value = "  +1.00%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.000009011"
$ws.Cells.Item(16, 5).Value = "  -0.49%  "
Set-TextValue $ws.Cells.Item(17, 4) "5.931"
$ws.Cells.Item(17, 5).Value = "  +0.75%  "
$ws.Cells.Item(18, 4).Value = "29.476.28"
$ws.Cells.Item(18, 5).Value = "  +1.71%  "
$ws.Cells.Item(19, 4).Value = "2.095.20"
$ws.Cells.Item(19, 5).Value = "  +0.50%  "
Set-TextValue $ws.Cells.Item(20, 4) "245.55"
$ws.Cells.Item(20, 5).Value = "  +2.68%  "
Set-TextValue $ws.Cells.Item(21, 4) "12.59"
$ws.Cells.Item(21, 5).Value = "  -0.67%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.9990"
$ws.Cells.Item(22, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(23, 4) "7.438"
$ws.Cells.Item(23, 5).Value = "  +3.44%  "
Set-TextValue $ws.Cells.Item(24, 4) "1.000"
$ws.Cells.Item(24, 5).Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(25, 4) "158.58"
$ws.Cells.Item(25, 5).Value = "  +0.22%  "
Set-TextValue $ws.Cells.Item(26, 4) "0.1418"
$ws.Cells.Item(26, 5).Value = "  +0.33%  "
Set-TextValue $ws.Cells.Item(27, 4) "8.598"
$ws.Cells.Item(27, 5).Value = "  +1.42%  "
Set-TextValue $ws.Cells.Item(28, 4) "17.83"
$ws.Cells.Item(28, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.501"
$ws.Cells.Item(29, 5).Value = "  +0.61%  "
Set-TextValue $ws.Cells.Item(30, 4) "0.05914"
$ws.Cells.Item(30, 5).Value = "  +5.76%  "
Set-TextValue $ws.Cells.Item(31, 4) "1.243"
$ws.Cells.Item(31, 5).Value = "  +2.58%  "
Set-TextValue $ws.Cells.Item(32, 4) "4.115"
$ws.Cells.Item(32, 5).Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(33, 4) "4.130"
$ws.Cells.Item(33, 5).Value = "  +0.24%  "
Set-TextValue $ws.Cells.Item(34, 4) "1.874"
$ws.Cells.Item(34, 5).Value = "  +1.85%  "
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(35, 4) "0.7306"
$ws.Cells.Item(35, 5).Value = "  -1.44%  "
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(36, 4) "1.145"
$ws.Cells.Item(36, 5).Value = "  +0.64%  "
Set-TextValue $ws.Cells.Item(37, 4) "2.612"
$ws.Cells.Item(37, 5).Value = "  -0.54%  "
Set-TextValue $ws.Cells.Item(38, 4) "2.887"
$ws.Cells.Item(38, 5).Value = "  +4.29%  "
$ws.Cells.Item(39, 4).Value = "1.233.27"
$ws.Cells.Item(39, 5).Value = "  +1.96%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.01766"
$ws.Cells.Item(40, 5).Value = "  -0.44%  "
$ws.Cells.Item(41, 2).Value = "FraxShare"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Cells.Item(41, 4) "6.288"
$ws.Cells.Item(41, 5).Value = "  -1.57%  "
$ws.Cells.Item(42, 2).Value = "TrustWalletToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Cells.Item(42, 4) "0.9175"
$ws.Cells.Item(42, 5).Value = "  +2.46%  "
$ws.Cells.Item(43, 5).Value = "  +0.44%  "
$ws.Cells.Item(44, 4).Value = "2.008.92"
$ws.Cells.Item(44, 5).Value = "  +1.12%  "
Set-TextValue $ws.Cells.Item(45, 4) "102.10"
$ws.Cells.Item(45, 5).Value = "  +0.87%  "
Set-TextValue $ws.Cells.Item(46, 4) "65.91"
$ws.Cells.Item(46, 5).Value = "  +1.04%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.00000000121"
$ws.Cells.Item(47, 5).Value = "  -0.94%  "
$ws.Cells.Item(48, 5).Value = "  -0.64%  "
$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(49, 4) "9.261"
$ws.Cells.Item(49, 5).Value = "  +1.51%  "
$ws.Cells.Item(50, 2).Value = "TheSandbox"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Cells.Item(50, 4) "0.4069"
$ws.Cells.Item(50, 5).Value = "  +0.41%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.1168"
$ws.Cells.Item(51, 5).Value = "  +5.99%  "
